{"js": "// Fix the typo \"thi\" -> \"this\" in \"...I prototyped them in python, thi should\n// be in the documents directory.\" and re-anchor the (invisible) \"_GoBack\"\n// bookmark next to the edit, mirroring where Word leaves it after the user\n// types the missing \"s\". The bookmark previously sat between \"w\" and \"hich\"\n// in the sentence above (\"...SOUT\u2026. which is connected to Bit 0)\"); that\n// split is gone once the bookmark moves, so that sentence simply reads\n// correctly as one continuous run of text.\n\nconst body = context.document.body;\n\n// 1) Drop the bookmark from its old (now stale) location. Word keeps exactly\n//    one \"_GoBack\" bookmark, always at the site of the most recent edit.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Insert the missing \"s\": find the unique \"thi\" that precedes \" should be\"\n//    and append \"s\" right after it (the search term purposely excludes the\n//    trailing space so InsertLocation.end lands exactly between \"thi\" and\n//    the space, not after it).\nlet results = body.search(\"python, thi\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n  target.insertText(\"s\", Word.InsertLocation.end);\n  await context.sync();\n}\n\n// 3) Re-create \"_GoBack\" as a collapsed bookmark right after the new \"s\",\n//    i.e. at the point Word just edited.\nresults = body.search(\"python, this\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n  const endPoint = target.getRange(Word.RangeLocation.end);\n  endPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Fix the typo \"thi\" -> \"this\" in \"...I prototyped them in python, thi should\n# be in the documents directory.\" and re-anchor the (invisible) \"_GoBack\"\n# bookmark next to the edit, mirroring where Word leaves it after the user\n# types the missing \"s\". The bookmark previously sat between \"w\" and \"hich\"\n# in the sentence above (\"...SOUT\u2026. which is connected to Bit 0)\"); once the\n# bookmark moves away from there, that sentence simply reads as one\n# continuous run of text.\n\n$d = $word.ActiveDocument\n\n# 1) Drop the bookmark from its old (now stale) location. Word keeps exactly\n#    one \"_GoBack\" bookmark, always at the site of the most recent edit.\n$oldMark = $d.Bookmarks(\"_GoBack\")\n$oldMark.Delete()\n\n# 2) Insert the missing \"s\": find the unique \"thi\" that precedes \" should be\"\n#    (the search term purposely excludes the trailing space) and drop the\n#    found range right after it, then insert \"s\" there.\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Execute(\"python, thi\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0) | Out-Null\n$range.Collapse(0)  # wdCollapseEnd\n$range.InsertAfter(\"s\")\n\n# 3) Re-create \"_GoBack\" as a collapsed bookmark right after the new \"s\",\n#    i.e. at the point Word just edited.\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.ClearFormatting()\n$find2.Execute(\"python, this\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0) | Out-Null\n$range2.Collapse(0)  # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $range2) | Out-Null\n"}
